# Error Calculations and Plots
# Applies the recorded edits to the missing_data worksheet:
#  - A few scattered cells toggle between a numeric value and "missing" (blank)
#  - Row 26 ("RM 232") is removed and row 28 ("SC 92") is removed, causing the
#    remaining rows (SC 5 .. SC 232) to shift up, with several of their values
#    corrected along the way. Net effect: the sheet shrinks from 35 data rows
#    (A1:F35) down to 33 (A1:F33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two trailing rows (34 and 35) so the sheet ends at row 33.
$ws.Rows.Item(34).Delete()
$ws.Rows.Item(34).Delete()

# Small "missing value" toggles scattered through rows 2-25 (row contents
# otherwise unchanged).
$ws.Range("E2").ClearContents()
$ws.Range("C6").Value = 15.1
$ws.Range("C8").ClearContents()
$ws.Range("C18").Value = 11.5
$ws.Range("C20").ClearContents()
$ws.Range("C23").Value = 12.2
$ws.Range("C25").ClearContents()

# Rows 26-33 now hold the corrected SC-series data (RM 232 was dropped and
# SC 92 was dropped, shifting everything else up and fixing some values).
$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = -14.6
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = 17

$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").ClearContents()
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

$ws.Range("A29").Value = "SC 119"
$ws.Range("B29").ClearContents()
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = -13
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06

$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18

$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").ClearContents()
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
